# Update "想去人数" (interested headcount) values for several events
# in both the "展览" sheet and the "全部类型" sheet.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 1368
$wsExhibit.Range("F3").Value = 2902
$wsExhibit.Range("F4").Value = 11

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 1368
$wsAll.Range("F4").Value = 2902
$wsAll.Range("F5").Value = 11
